$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at 401-405; this shifts the existing rows 401-411 down to 406-416.
$ws.Range("A401:A405").EntireRow.Insert()

# Common (constant) column values shared by every data row in this block.
$commonA = 3
$commonB = "Femacal de La Calera"
$commonC = "Coquimbo"
$commonE = 5
$commonF = "Fruta"
$commonG = 100102
$commonH = "Cítricos"
$commonI = 100102005
$commonJ = "Naranja"
$commonQ = "`$/malla 13 kilos"
$commonR = "Provincia de Quillota"
$commonT = 13

# New row data: row => @(D,K,L,M,N,O,P,S)
$newRows = @{
    401 = @(44448, "Cara cara",  "Primera", 70,  4000, 4000, 4000, 308)
    402 = @(44448, "Lane Late",  "Primera", 186, 3500, 4000, 3739, 288)
    403 = @(44448, "Lane Late",  "Segunda", 177, 2500, 3000, 2754, 212)
    404 = @(44448, "Navel Late", "Primera", 172, 3500, 3800, 3648, 281)
    405 = @(44448, "Navel Late", "Segunda", 140, 2500, 3000, 2786, 214)
}

foreach ($r in 401..405) {
    $vals = $newRows[$r]

    $ws.Cells.Item($r, 1).Value = $commonA
    $ws.Cells.Item($r, 2).Value = $commonB
    $ws.Cells.Item($r, 3).Value = $commonC
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $commonE
    $ws.Cells.Item($r, 6).Value = $commonF
    $ws.Cells.Item($r, 7).Value = $commonG
    $ws.Cells.Item($r, 8).Value = $commonH
    $ws.Cells.Item($r, 9).Value = $commonI
    $ws.Cells.Item($r, 10).Value = $commonJ
    $ws.Cells.Item($r, 11).Value = $vals[1]
    $ws.Cells.Item($r, 12).Value = $vals[2]
    $ws.Cells.Item($r, 13).Value = $vals[3]
    $ws.Cells.Item($r, 14).Value = $vals[4]
    $ws.Cells.Item($r, 15).Value = $vals[5]
    $ws.Cells.Item($r, 16).Value = $vals[6]
    $ws.Cells.Item($r, 17).Value = $commonQ
    $ws.Cells.Item($r, 18).Value = $commonR
    $ws.Cells.Item($r, 19).Value = $vals[7]
    $ws.Cells.Item($r, 20).Value = $commonT
}
